$wb = $excel.ActiveWorkbook

# ---- Sheet "ROW50-FE-LIFTER" -> append new row 39 ----
$ws = $wb.Worksheets.Item("ROW50-FE-LIFTER")
$r = 39
$ws.Cells.Item($r, 1).Value = 45743.66920777778
$ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($r, 2).Value = "0x01,0x90"
$ws.Cells.Item($r, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
$ws.Cells.Item($r, 4).Value = "0x01,0x6e"
$ws.Cells.Item($r, 5).Value = "0xe"
$ws.Cells.Item($r, 6).Value = 400
$ws.Cells.Item($r, 7).Value = [double]"5.68631262647114e+23"
$ws.Cells.Item($r, 8).Value = 366
$ws.Cells.Item($r, 9).Value = 14

# ---- Sheet "ROW50-MID-LIFTER" -> append new row 41 ----
$ws = $wb.Worksheets.Item("ROW50-MID-LIFTER")
$r = 41
$ws.Cells.Item($r, 1).Value = 45743.63802083334
$ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($r, 2).Value = "0x01,0x90 "
$ws.Cells.Item($r, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws.Cells.Item($r, 4).Value = "0x01,0x72"
$ws.Cells.Item($r, 5).Value = "0x19"
$ws.Cells.Item($r, 6).Value = 400
$ws.Cells.Item($r, 7).Value = "'568631262647113771663628"
$ws.Cells.Item($r, 7).Style = "Normal"
$ws.Cells.Item($r, 8).Value = 370
$ws.Cells.Item($r, 9).Value = 25

# ---- Sheet "ROW11-FE-LIFTER" -> append new row 39 ----
$ws = $wb.Worksheets.Item("ROW11-FE-LIFTER")
$r = 39
$ws.Cells.Item($r, 1).Value = 45743.68740349537
$ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($r, 2).Value = "0x01,0x90"
$ws.Cells.Item($r, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
$ws.Cells.Item($r, 4).Value = "0x01,0x6e"
$ws.Cells.Item($r, 5).Value = "0x14"
$ws.Cells.Item($r, 6).Value = 400
$ws.Cells.Item($r, 7).Value = [double]"5.68631262647114e+23"
$ws.Cells.Item($r, 8).Value = 366
$ws.Cells.Item($r, 9).Value = 20

# ---- Sheet "ROW11-MID-LIFTER" -> append new row 39 ----
$ws = $wb.Worksheets.Item("ROW11-MID-LIFTER")
$r = 39
$ws.Cells.Item($r, 1).Value = 45743.83452253472
$ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($r, 2).Value = "0x01,0x90"
$ws.Cells.Item($r, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
$ws.Cells.Item($r, 4).Value = "0x01,0x76"
$ws.Cells.Item($r, 5).Value = "0x19"
$ws.Cells.Item($r, 6).Value = 400
$ws.Cells.Item($r, 7).Value = [double]"5.68631262647114e+23"
$ws.Cells.Item($r, 8).Value = 374
$ws.Cells.Item($r, 9).Value = 25
